$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = ""

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "أحمد شريم"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "2"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "الصمود"

$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "الرحلة 2"

$ws.Cells.Item($row, 6).NumberFormat = "@"
$ws.Cells.Item($row, 6).Value = "C2"

$ws.Cells.Item($row, 7).NumberFormat = "@"
$ws.Cells.Item($row, 7).Value = "NRC"

$ws.Cells.Item($row, 8).NumberFormat = "@"
$ws.Cells.Item($row, 8).Value = "٠٦:٠١:٥٥ م"
